$d = $word.ActiveDocument

# --- List 2 style: firstLine 475 -> left 720 / hanging 360 ---------------
$list2 = $d.Styles("List 2")
$list2.ParagraphFormat.LeftIndent = 36      # 720 twips = 36 pt
$list2.ParagraphFormat.FirstLineIndent = -18  # hanging 360 twips = 18 pt

# --- List 1 style: drop the firstLine indent ------------------------------
$list1 = $d.Styles("List 1")
$list1.ParagraphFormat.FirstLineIndent = 0

# --- List 3 (new, built-in) style: left 1080 / hanging 360, contextual spacing
# Built-in styles only get minted into styles.xml once referenced by a
# paragraph, so apply it momentarily and then restore the paragraph's
# original style - this leaves the style definition behind without
# altering any visible document content.
$p = $d.Paragraphs.First
$savedStyle = $p.Range.Style
$p.Range.Style = $d.Styles(-52)          # wdStyleList3 ("List 3")

$list3 = $d.Styles("List 3")
$list3.ParagraphFormat.LeftIndent = 54      # 1080 twips = 54 pt
$list3.ParagraphFormat.FirstLineIndent = -18  # hanging 360 twips = 18 pt
$list3.NoSpaceBetweenParagraphsOfSameStyle = $true

$p.Range.Style = $savedStyle
